$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Labadie"
$ws.Range("B4").Value = "Marian Altenwerth I"
$ws.Range("B6").Value = "Actualizado"
$ws.Range("B7").Value = "Actualizado"

$ws.Range("B7").Select()
